# Replace the 25 three-digit-by-one-digit multiplication equations in the
# practice table with a new set of problems (commit "Update master to
# output generated at 503736d"). Each old equation string is unique in the
# document, so a plain exact-text Find/Replace is safe for every cell.
$d = $word.ActiveDocument

$d.Content.Find.Execute("153×7=1071", $true, $false, $false, $false, $false, $true, 1, $false, "120×5=600", 2) | Out-Null
$d.Content.Find.Execute("506×9=4554", $true, $false, $false, $false, $false, $true, 1, $false, "345×9=3105", 2) | Out-Null
$d.Content.Find.Execute("993×7=6951", $true, $false, $false, $false, $false, $true, 1, $false, "316×3=948", 2) | Out-Null
$d.Content.Find.Execute("124×6=744", $true, $false, $false, $false, $false, $true, 1, $false, "219×8=1752", 2) | Out-Null
$d.Content.Find.Execute("220×4=880", $true, $false, $false, $false, $false, $true, 1, $false, "654×7=4578", 2) | Out-Null
$d.Content.Find.Execute("517×8=4136", $true, $false, $false, $false, $false, $true, 1, $false, "564×9=5076", 2) | Out-Null
$d.Content.Find.Execute("235×4=940", $true, $false, $false, $false, $false, $true, 1, $false, "491×9=4419", 2) | Out-Null
$d.Content.Find.Execute("239×9=2151", $true, $false, $false, $false, $false, $true, 1, $false, "207×5=1035", 2) | Out-Null
$d.Content.Find.Execute("668×4=2672", $true, $false, $false, $false, $false, $true, 1, $false, "525×3=1575", 2) | Out-Null
$d.Content.Find.Execute("758×7=5306", $true, $false, $false, $false, $false, $true, 1, $false, "115×5=575", 2) | Out-Null
$d.Content.Find.Execute("800×4=3200", $true, $false, $false, $false, $false, $true, 1, $false, "896×5=4480", 2) | Out-Null
$d.Content.Find.Execute("141×7=987", $true, $false, $false, $false, $false, $true, 1, $false, "839×9=7551", 2) | Out-Null
$d.Content.Find.Execute("971×2=1942", $true, $false, $false, $false, $false, $true, 1, $false, "184×4=736", 2) | Out-Null
$d.Content.Find.Execute("763×9=6867", $true, $false, $false, $false, $false, $true, 1, $false, "442×9=3978", 2) | Out-Null
$d.Content.Find.Execute("741×5=3705", $true, $false, $false, $false, $false, $true, 1, $false, "332×4=1328", 2) | Out-Null
$d.Content.Find.Execute("364×8=2912", $true, $false, $false, $false, $false, $true, 1, $false, "249×3=747", 2) | Out-Null
$d.Content.Find.Execute("852×6=5112", $true, $false, $false, $false, $false, $true, 1, $false, "361×5=1805", 2) | Out-Null
$d.Content.Find.Execute("253×7=1771", $true, $false, $false, $false, $false, $true, 1, $false, "692×5=3460", 2) | Out-Null
$d.Content.Find.Execute("991×5=4955", $true, $false, $false, $false, $false, $true, 1, $false, "911×8=7288", 2) | Out-Null
$d.Content.Find.Execute("532×9=4788", $true, $false, $false, $false, $false, $true, 1, $false, "801×4=3204", 2) | Out-Null
$d.Content.Find.Execute("193×3=579", $true, $false, $false, $false, $false, $true, 1, $false, "782×2=1564", 2) | Out-Null
$d.Content.Find.Execute("493×7=3451", $true, $false, $false, $false, $false, $true, 1, $false, "165×3=495", 2) | Out-Null
$d.Content.Find.Execute("808×5=4040", $true, $false, $false, $false, $false, $true, 1, $false, "351×9=3159", 2) | Out-Null
$d.Content.Find.Execute("610×8=4880", $true, $false, $false, $false, $false, $true, 1, $false, "264×6=1584", 2) | Out-Null
$d.Content.Find.Execute("640×3=1920", $true, $false, $false, $false, $false, $true, 1, $false, "101×6=606", 2) | Out-Null
